$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number format to A2:A6 and C2:C6 (REF_NO / TRN_REF_NO style columns)
$ws.Range("A2:A6").NumberFormat = "@"
$ws.Range("C2:C6").NumberFormat = "@"

# Apply text number format to B2:B6 and replace the numeric bank-code values
# with their text equivalents (so they become shared strings on save)
$ws.Range("B2:B6").NumberFormat = "@"
$ws.Range("B2").Value = "199999"
$ws.Range("B3").Value = "288888"
$ws.Range("B4").Value = "388888"
$ws.Range("B5").Value = "488888"
$ws.Range("B6").Value = "588888"

# Apply a 2-decimal numeric format to the amount column
$ws.Range("E2:E6").NumberFormat = "0.00"

# Move the active selection to C3
$ws.Range("C3").Select()

# Switch the page to portrait orientation
$ws.PageSetup.Orientation = 1
